$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D15 and D17 venue values: "bioRxiv preprint" -> "preprint" (matching D14/D16)
$ws.Range("D15").Value = "preprint"
$ws.Range("D17").Value = "preprint"

# Set column J width
$ws.Columns.Item(10).ColumnWidth = 21.5

# Update selection to D18
$ws.Range("D18").Select()
